$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function Find-ParagraphIndex($Contains) {
    $d = $word.ActiveDocument
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        if ($t -like "*$Contains*") {
            return $i
        }
    }
    return -1
}

function Set-ParagraphXml($Index, $InnerXml) {
    $d = $word.ActiveDocument
    $p = $d.Paragraphs.Item($Index)
    $r = $p.Range
    $r.InsertXML($InnerXml)
}

# --- 1) Merge "-1-  " / "Figure out what we need to do with the dashboard" /
#        " (priority High)" into a single run (keeping the leading tab run). ---
$idx1 = Find-ParagraphIndex "-1-  Figure out what we need to do with the dashboard (priority High)"
$xml1 = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="14922B17" w14:textId="7219A1D5" w:rsidR="00B25679" w:rsidRDefault="00B25679" w:rsidP="00B25679"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/><w:t>-1-  Figure out what we need to do with the dashboard (priority High)</w:t></w:r></w:p>'
Set-ParagraphXml $idx1 $xml1

# --- 2) Merge "-2- " / "check how to connect the pipeline to the dashboard" /
#        " (priority High)" into a single run. ---
$idx2 = Find-ParagraphIndex "-2- check how to connect the pipeline to the dashboard (priority High)"
$xml2 = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="4A44226A" w14:textId="1FE3DFDF" w:rsidR="00B25679" w:rsidRDefault="00B25679" w:rsidP="00B25679"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/><w:t>-2- check how to connect the pipeline to the dashboard (priority High)</w:t></w:r></w:p>'
Set-ParagraphXml $idx2 $xml2

# --- 3) Merge "-5- " / "search the simplest stack for the dashboard" /
#        " (meduim)" into a single run. ---
$idx3 = Find-ParagraphIndex "-5- search the simplest stack for the dashboard (meduim)"
$xml3 = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="796CB010" w14:textId="7C4E34F7" w:rsidR="00B25679" w:rsidRDefault="00B25679" w:rsidP="00B25679"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/><w:t>-5- search the simplest stack for the dashboard (meduim)</w:t></w:r></w:p>'
Set-ParagraphXml $idx3 $xml3

# --- 4) "Search the simplest stack for the dashboard ()" -> split the run so
#        "Aadit" is inserted between the parentheses, as 3 separate runs. ---
$idx4 = Find-ParagraphIndex "dashboard ()"
$xml4 = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="42B64C05" w14:textId="634BE8BA" w:rsidR="007F5A5E" w:rsidRDefault="007F5A5E" w:rsidP="00503EF8"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/><w:t>Search the simplest stack for the dashboard (</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Aadit</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>)</w:t></w:r></w:p>'
Set-ParagraphXml $idx4 $xml4

Write-Output "Edits applied: idx1=$idx1 idx2=$idx2 idx3=$idx3 idx4=$idx4"
